# Scheduled market-data refresh: update crafting leve profit figures
# (currentAveragePrice / NQ / HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets with freshly
# pulled market-board values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 246.70589
$ws.Range("I53").Value = 239.09091
$ws.Range("J53").Value = 260.66666
$ws.Range("K53").Value = 239.09091
$ws.Range("L53").Value = 260.66666
$ws.Range("M53").Value = 397.90909
$ws.Range("N53").Value = -1534.66666

$ws.Range("H62").Value = 5526.154
$ws.Range("I62").Value = 3980
$ws.Range("K62").Value = 3980
$ws.Range("M62").Value = -3356

$ws.Range("H64").Value = 3194
$ws.Range("I64").Value = 3180
$ws.Range("J64").Value = 3197.5
$ws.Range("K64").Value = 3180
$ws.Range("L64").Value = 3197.5
$ws.Range("M64").Value = -2932
$ws.Range("N64").Value = -3693.5

$ws.Range("H65").Value = 5526.154
$ws.Range("I65").Value = 3980
$ws.Range("K65").Value = 19900
$ws.Range("M65").Value = -16780

$ws.Range("H67").Value = 3194
$ws.Range("I67").Value = 3180
$ws.Range("J67").Value = 3197.5
$ws.Range("K67").Value = 3180
$ws.Range("L67").Value = 3197.5
$ws.Range("M67").Value = -2322
$ws.Range("N67").Value = -4913.5

$ws.Range("H69").Value = 7625

$ws.Range("H72").Value = 7625

$ws.Range("H132").Value = 1603.5807
$ws.Range("J132").Value = 1698
$ws.Range("L132").Value = 5094
$ws.Range("N132").Value = -10154

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = ""

$ws.Range("H135").Value = 812.64
$ws.Range("I135").Value = 537.2222
$ws.Range("J135").Value = 1520.8572
$ws.Range("K135").Value = 4834.999800000001
$ws.Range("L135").Value = 13687.7148
$ws.Range("M135").Value = -2299.999800000001
$ws.Range("N135").Value = -18757.7148

$ws.Range("H137").Value = 2961.1785
$ws.Range("I137").Value = 1966.5714
$ws.Range("J137").Value = 3955.7856
$ws.Range("K137").Value = 5899.7142
$ws.Range("L137").Value = 11867.3568
$ws.Range("M137").Value = -3349.7142
$ws.Range("N137").Value = -16967.3568

$ws.Range("H138").Value = 13181
$ws.Range("I138").Value = 11756.2
$ws.Range("J138").Value = 14071.5
$ws.Range("K138").Value = 35268.60000000001
$ws.Range("L138").Value = 42214.5
$ws.Range("M138").Value = -30128.60000000001
$ws.Range("N138").Value = -52494.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14613.794
$ws.Range("I32").Value = 6216.6113
$ws.Range("J32").Value = 25810.037
$ws.Range("K32").Value = 6216.6113
$ws.Range("L32").Value = 25810.037
$ws.Range("M32").Value = -5929.6113
$ws.Range("N32").Value = -26384.037

$ws.Range("H88").Value = 2497.375
$ws.Range("I88").Value = 1997.4
$ws.Range("J88").Value = 3330.6667
$ws.Range("K88").Value = 1997.4
$ws.Range("L88").Value = 3330.6667
$ws.Range("M88").Value = -1591.4
$ws.Range("N88").Value = -4142.6667

$ws.Range("H91").Value = 2497.375
$ws.Range("I91").Value = 1997.4
$ws.Range("J91").Value = 3330.6667
$ws.Range("K91").Value = 1997.4
$ws.Range("L91").Value = 3330.6667
$ws.Range("M91").Value = -593.4000000000001
$ws.Range("N91").Value = -6138.6667

$ws.Range("H97").Value = 648.7778
$ws.Range("I97").Value = 648.7778
$ws.Range("K97").Value = 648.7778
$ws.Range("M97").Value = -152.7778

$ws.Range("H102").Value = 2220.111
$ws.Range("I102").Value = 1955.6666
$ws.Range("J102").Value = 2749
$ws.Range("K102").Value = 1955.6666
$ws.Range("L102").Value = 2749
$ws.Range("M102").Value = -333.6666
$ws.Range("N102").Value = -5993

$ws.Range("H122").Value = 419704.75
$ws.Range("I122").Value = 668809.3
$ws.Range("K122").Value = 2006427.9
$ws.Range("M122").Value = -2003977.9

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""

$ws.Range("H132").Value = 1847.84
$ws.Range("I132").Value = 1803.9183
$ws.Range("K132").Value = 5411.7549
$ws.Range("M132").Value = -2881.7549

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1735.5264
$ws.Range("I99").Value = 1540
$ws.Range("J99").Value = 2070.7144
$ws.Range("K99").Value = 1540
$ws.Range("L99").Value = 2070.7144
$ws.Range("M99").Value = -42
$ws.Range("N99").Value = -5066.7144

$ws.Range("H130").Value = 54498.5
$ws.Range("J130").Value = 54498.5
$ws.Range("L130").Value = 54498.5
$ws.Range("N130").Value = -64538.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4434.8
$ws.Range("I31").Value = 2144.3333
$ws.Range("J31").Value = 5416.4287
$ws.Range("K31").Value = 2144.3333
$ws.Range("L31").Value = 5416.4287
$ws.Range("M31").Value = -1849.3333
$ws.Range("N31").Value = -6006.4287

$ws.Range("H34").Value = 4434.8
$ws.Range("I34").Value = 2144.3333
$ws.Range("J34").Value = 5416.4287
$ws.Range("K34").Value = 2144.3333
$ws.Range("L34").Value = 5416.4287
$ws.Range("M34").Value = -1942.3333
$ws.Range("N34").Value = -5820.4287

$ws.Range("H58").Value = 6378
$ws.Range("I58").Value = 3498
$ws.Range("J58").Value = 6954
$ws.Range("K58").Value = 3498
$ws.Range("L58").Value = 6954
$ws.Range("M58").Value = -3295
$ws.Range("N58").Value = -7360

$ws.Range("H107").Value = 570.1875
$ws.Range("I107").Value = 461.66666
$ws.Range("J107").Value = 709.7143
$ws.Range("K107").Value = 461.66666
$ws.Range("L107").Value = 709.7143
$ws.Range("M107").Value = 1458.33334
$ws.Range("N107").Value = -4549.7143

$ws.Range("H132").Value = 3327.3215
$ws.Range("I132").Value = 2962.762
$ws.Range("J132").Value = 4421
$ws.Range("K132").Value = 8888.286
$ws.Range("L132").Value = 13263
$ws.Range("M132").Value = -6358.286
$ws.Range("N132").Value = -18323

$ws.Range("H134").Value = 3919.15
$ws.Range("I134").Value = 2724.4666
$ws.Range("K134").Value = 8173.399800000001
$ws.Range("M134").Value = -5638.399800000001

$ws.Range("H136").Value = 6378
$ws.Range("I136").Value = 3498
$ws.Range("J136").Value = 6954
$ws.Range("K136").Value = 10494
$ws.Range("L136").Value = 20862
$ws.Range("M136").Value = -7944
$ws.Range("N136").Value = -25962

$ws.Range("H138").Value = 29999
$ws.Range("J138").Value = 29999
$ws.Range("L138").Value = 29999
$ws.Range("N138").Value = -40279

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2282.842
$ws.Range("I132").Value = 1058.3334
$ws.Range("K132").Value = 9525.000599999999
$ws.Range("M132").Value = -6995.000599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 18000
$ws.Range("I22").Value = 18000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 18000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -17471
$ws.Range("N22").Value = ""

$ws.Range("H80").Value = 4033.6667
$ws.Range("I80").Value = 3249.5
$ws.Range("J80").Value = 4425.75
$ws.Range("K80").Value = 3249.5
$ws.Range("L80").Value = 4425.75
$ws.Range("M80").Value = -2251.5
$ws.Range("N80").Value = -6421.75

$ws.Range("H83").Value = 4033.6667
$ws.Range("I83").Value = 3249.5
$ws.Range("J83").Value = 4425.75
$ws.Range("K83").Value = 16247.5
$ws.Range("L83").Value = 22128.75
$ws.Range("M83").Value = -11255.5
$ws.Range("N83").Value = -32112.75

$ws.Range("H97").Value = 1897.1072
$ws.Range("I97").Value = 2032.6
$ws.Range("K97").Value = 2032.6
$ws.Range("M97").Value = -1536.6

$ws.Range("H113").Value = 3426.4375
$ws.Range("J113").Value = 3601
$ws.Range("L113").Value = 3601
$ws.Range("N113").Value = -7941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 20000000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = ""

$ws.Range("H11").Value = 5000
$ws.Range("J11").Value = 5000
$ws.Range("L11").Value = 5000
$ws.Range("N11").Value = -5280

$ws.Range("H22").Value = 937.0769
$ws.Range("I22").Value = 1050
$ws.Range("K22").Value = 1050
$ws.Range("M22").Value = -755

$ws.Range("H27").Value = 937.0769
$ws.Range("I27").Value = 1050
$ws.Range("K27").Value = 1050
$ws.Range("M27").Value = -943

$ws.Range("H43").Value = 20000
$ws.Range("J43").Value = 20000
$ws.Range("L43").Value = 20000
$ws.Range("N43").Value = -20386

$ws.Range("H46").Value = 2781.3125
$ws.Range("I46").Value = 1112.625
$ws.Range("J46").Value = 4450
$ws.Range("K46").Value = 1112.625
$ws.Range("L46").Value = 4450
$ws.Range("M46").Value = -924.625
$ws.Range("N46").Value = -4826

$ws.Range("H93").Value = 845.0454999999999
$ws.Range("I93").Value = 635.6429000000001
$ws.Range("K93").Value = 635.6429000000001
$ws.Range("M93").Value = 612.3570999999999

$ws.Range("H132").Value = 4419.815
$ws.Range("I132").Value = 3120.1538
$ws.Range("K132").Value = 9360.4614
$ws.Range("M132").Value = -6830.4614

$ws.Range("H136").Value = 4042.111
$ws.Range("I136").Value = 3922.375
$ws.Range("K136").Value = 11767.125
$ws.Range("M136").Value = -9217.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 653.8570999999999
$ws.Range("I96").Value = 634.2
$ws.Range("J96").Value = 703
$ws.Range("K96").Value = 634.2
$ws.Range("L96").Value = 703
$ws.Range("M96").Value = 738.8
$ws.Range("N96").Value = -3449

$ws.Range("H132").Value = 2899
$ws.Range("I132").Value = 2348.5
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 7045.5
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -4515.5
$ws.Range("N132").Value = -17060
